$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 225, pushing existing rows 225-256 down to 227-258.
$ws.Rows.Item(225).Resize(2).Insert()

# New row 225 data
$ws.Cells.Item(225, 1).Value = 9
$ws.Cells.Item(225, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(225, 3).Value = "Metropolitana"
$ws.Cells.Item(225, 4).Value = 44474
$ws.Cells.Item(225, 4).NumberFormat = $ws.Cells.Item(224, 4).NumberFormat
$ws.Cells.Item(225, 5).Value = 13
$ws.Cells.Item(225, 6).Value = 100112031
$ws.Cells.Item(225, 7).Value = "Poroto verde"
$ws.Cells.Item(225, 8).Value = "Magnum"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 61
$ws.Cells.Item(225, 11).Value = 43000
$ws.Cells.Item(225, 12).Value = 45000
$ws.Cells.Item(225, 13).Value = 43984
$ws.Cells.Item(225, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(225, 15).Value = "Perú"
$ws.Cells.Item(225, 16).Value = 1759
$ws.Cells.Item(225, 17).Value = 25
$ws.Cells.Item(225, 18).Value = "Hortaliza"

# New row 226 data
$ws.Cells.Item(226, 1).Value = 9
$ws.Cells.Item(226, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(226, 3).Value = "Metropolitana"
$ws.Cells.Item(226, 4).Value = 44474
$ws.Cells.Item(226, 4).NumberFormat = $ws.Cells.Item(224, 4).NumberFormat
$ws.Cells.Item(226, 5).Value = 13
$ws.Cells.Item(226, 6).Value = 100112031
$ws.Cells.Item(226, 7).Value = "Poroto verde"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 43
$ws.Cells.Item(226, 11).Value = 38000
$ws.Cells.Item(226, 12).Value = 39000
$ws.Cells.Item(226, 13).Value = 38512
$ws.Cells.Item(226, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(226, 15).Value = "Perú"
$ws.Cells.Item(226, 16).Value = 1540
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
